# Applies the scheduled-runner data refresh to the Kujata_Profits workbook.
# Updates computed market-price / profit columns (H-N) across all eight job sheets
# to match the latest Universalis price snapshot; one stale cell (CRP!M22) is cleared
# because that leve no longer has an HQ-profit figure in the refreshed dataset.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 883629.1
$ws.Range("I88").Value = 1507.4
$ws.Range("J88").Value = 1373696.8
$ws.Range("K88").Value = 1507.4
$ws.Range("L88").Value = 1373696.8
$ws.Range("M88").Value = -1101.4
$ws.Range("N88").Value = -1374508.8
$ws.Range("H91").Value = 883629.1
$ws.Range("I91").Value = 1507.4
$ws.Range("J91").Value = 1373696.8
$ws.Range("K91").Value = 1507.4
$ws.Range("L91").Value = 1373696.8
$ws.Range("M91").Value = -103.4000000000001
$ws.Range("N91").Value = -1376504.8
$ws.Range("H133").Value = 35000
$ws.Range("J133").Value = 35000
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -45120
$ws.Range("H137").Value = 1278.5385
$ws.Range("I137").Value = 886.96875
$ws.Range("K137").Value = 2660.90625
$ws.Range("M137").Value = -110.90625
$ws.Range("H138").Value = 1206.2988
$ws.Range("I138").Value = 720.70734
$ws.Range("J138").Value = 1639.1086
$ws.Range("K138").Value = 2162.12202
$ws.Range("L138").Value = 4917.325800000001
$ws.Range("M138").Value = 2977.87798
$ws.Range("N138").Value = -15197.3258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1578.4286
$ws.Range("I74").Value = 1125.4706
$ws.Range("J74").Value = 3503.5
$ws.Range("K74").Value = 1125.4706
$ws.Range("L74").Value = 3503.5
$ws.Range("M74").Value = -251.4706000000001
$ws.Range("N74").Value = -5251.5
$ws.Range("H77").Value = 1578.4286
$ws.Range("I77").Value = 1125.4706
$ws.Range("J77").Value = 3503.5
$ws.Range("K77").Value = 5627.353000000001
$ws.Range("L77").Value = 17517.5
$ws.Range("M77").Value = -1259.353000000001
$ws.Range("N77").Value = -26253.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1329.6086
$ws.Range("I107").Value = 1012.7059
$ws.Range("K107").Value = 1012.7059
$ws.Range("M107").Value = 907.2941
$ws.Range("H134").Value = 4507.303
$ws.Range("I134").Value = 1183.6897
$ws.Range("K134").Value = 3551.0691
$ws.Range("M134").Value = -1016.0691

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 700000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H31").Value = 2233.4092
$ws.Range("I31").Value = 2246.85
$ws.Range("J31").Value = 2099
$ws.Range("K31").Value = 2246.85
$ws.Range("L31").Value = 2099
$ws.Range("M31").Value = -1951.85
$ws.Range("N31").Value = -2689
$ws.Range("H34").Value = 2233.4092
$ws.Range("I34").Value = 2246.85
$ws.Range("J34").Value = 2099
$ws.Range("K34").Value = 2246.85
$ws.Range("L34").Value = 2099
$ws.Range("M34").Value = -2044.85
$ws.Range("N34").Value = -2503
$ws.Range("H62").Value = 5716597
$ws.Range("I62").Value = 2379.3823
$ws.Range("J62").Value = 200000000
$ws.Range("K62").Value = 2379.3823
$ws.Range("L62").Value = 200000000
$ws.Range("M62").Value = -1755.3823
$ws.Range("N62").Value = -200001248
$ws.Range("H65").Value = 5716597
$ws.Range("I65").Value = 2379.3823
$ws.Range("J65").Value = 200000000
$ws.Range("K65").Value = 11896.9115
$ws.Range("L65").Value = 1000000000
$ws.Range("M65").Value = -8776.911500000002
$ws.Range("N65").Value = -1000006240
$ws.Range("H107").Value = 559.381
$ws.Range("I107").Value = 479.08334
$ws.Range("J107").Value = 666.44446
$ws.Range("K107").Value = 479.08334
$ws.Range("L107").Value = 666.44446
$ws.Range("M107").Value = 1440.91666
$ws.Range("N107").Value = -4506.44446
$ws.Range("H132").Value = 3214.7856
$ws.Range("I132").Value = 3432.28
$ws.Range("J132").Value = 1402.3334
$ws.Range("K132").Value = 10296.84
$ws.Range("L132").Value = 4207.0002
$ws.Range("M132").Value = -7766.84
$ws.Range("N132").Value = -9267.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 599.1111
$ws.Range("I5").Value = 570.2857
$ws.Range("K5").Value = 1710.8571
$ws.Range("M5").Value = -1598.8571
$ws.Range("H11").Value = 230.92308
$ws.Range("I11").Value = 276.83334
$ws.Range("J11").Value = 191.57143
$ws.Range("K11").Value = 830.5000200000001
$ws.Range("L11").Value = 574.71429
$ws.Range("M11").Value = -690.5000200000001
$ws.Range("N11").Value = -854.71429
$ws.Range("H107").Value = 5111.857
$ws.Range("J107").Value = 7496.9287
$ws.Range("L107").Value = 22490.7861
$ws.Range("N107").Value = -26330.7861
$ws.Range("H122").Value = 754.0645
$ws.Range("I122").Value = 628
$ws.Range("J122").Value = 845.1111
$ws.Range("K122").Value = 5652
$ws.Range("L122").Value = 7605.9999
$ws.Range("M122").Value = -3202
$ws.Range("N122").Value = -12505.9999
$ws.Range("H123").Value = 2170.2666
$ws.Range("I123").Value = 1080
$ws.Range("J123").Value = 2897.111
$ws.Range("K123").Value = 3240
$ws.Range("L123").Value = 8691.332999999999
$ws.Range("M123").Value = -790
$ws.Range("N123").Value = -13591.333
$ws.Range("H135").Value = 599.1111
$ws.Range("I135").Value = 570.2857
$ws.Range("K135").Value = 5132.571300000001
$ws.Range("M135").Value = -2597.571300000001
$ws.Range("H136").Value = 3419
$ws.Range("I136").Value = 2388.3333
$ws.Range("J136").Value = 4449.6665
$ws.Range("K136").Value = 7164.999899999999
$ws.Range("L136").Value = 13348.9995
$ws.Range("M136").Value = -2064.999899999999
$ws.Range("N136").Value = -23548.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1563.1666
$ws.Range("I113").Value = 1480.625
$ws.Range("J113").Value = 1728.25
$ws.Range("K113").Value = 1480.625
$ws.Range("L113").Value = 1728.25
$ws.Range("M113").Value = 689.375
$ws.Range("N113").Value = -6068.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1527.6
$ws.Range("I61").Value = 1308.25
$ws.Range("K61").Value = 1308.25
$ws.Range("M61").Value = -1106.25
$ws.Range("H68").Value = 1936
$ws.Range("I68").Value = 1920
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1920
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1171
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 1936
$ws.Range("I71").Value = 1920
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 9600
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -5856
$ws.Range("N71").Value = -17488
$ws.Range("H113").Value = 1527.6
$ws.Range("I113").Value = 1308.25
$ws.Range("K113").Value = 1308.25
$ws.Range("M113").Value = 861.75
$ws.Range("H127").Value = 36666.668
$ws.Range("J127").Value = 36666.668
$ws.Range("L127").Value = 36666.668
$ws.Range("N127").Value = -46586.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 553
$ws.Range("I107").Value = 476.2857
$ws.Range("K107").Value = 1428.8571
$ws.Range("M107").Value = 491.1428999999998
$ws.Range("H132").Value = 2731.625
$ws.Range("I132").Value = 3334.0588
$ws.Range("J132").Value = 1268.5714
$ws.Range("K132").Value = 10002.1764
$ws.Range("L132").Value = 3805.7142
$ws.Range("M132").Value = -7472.1764
$ws.Range("N132").Value = -8865.7142
